$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = "Κείμενο άσκησης" (exercise text/question), Column E = "Λύση άσκησης" (solution).
# Set E10 first so the shared-string table assigns it the lower index (84),
# leaving D10's (new) string at index 85 - matching the target workbook's
# shared-string ordering.
$e10Text = @"
Για τον υπολογιστή **172.16.150.10/20** να υπολογίσετε:  <br>
**Δ1**. Την μάσκα δικτύου(δυαδική-δεκαδική)<br>
255.255.240.0 ή 11111111.11111111.1111 **0000.00000000**<br>
**Δ2.** Τη διεύθυνση δικτύου (network address)  <br>
Διεύθυνση Δικτύου   <br>
172.16.150.10(1010 **0110.00001010**) AND 255.255.240.0(1111 **00000.00000000**)=172.16.144.0(1001 **0000.00000000**)/20 <br>
**Δ3.** Τη διεύθυνση εκπομπής (broadcast address)  <br>
Διεύθυνση Εκπομπής ->  172.16.159.255(1001 **1111.11111111**)  <br>
**Δ4.** Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου  <br>
2^12 -2 = 4094 
**Δ5.** Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή  \
Από 172.16.144.1 έως 172.16.159.254  \
"@

$d10Text = @"
Για τον υπολογιστή <strong>172.16.150.10/20</strong> να υπολογίσετε:  <br>
**Δ1.** Την μάσκα δικτύου(δυαδική-δεκαδική)  <br>
**Δ2.** Τη διεύθυνση δικτύου (network address)  <br> 
**Δ3.** Τη διεύθυνση εκπομπής (broadcast address)  <br>
**Δ4.** Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου  <br>
**Δ5.** Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή  <br>
"@

$ws.Range("E10").Value = $e10Text
$ws.Range("D10").Value = $d10Text

# Editing the wrapped cells auto-expands row 10; restore its original height.
$ws.Rows.Item(10).RowHeight = 14.25

$null = $ws.Range("D11").Select()
